$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round Q2 and R2 to nearest integer
$ws.Range("Q2").Value = 723645
$ws.Range("R2").Value = 7330384

# Remove Z2 (Starttid) and AB2 (Sluttid) cell contents entirely
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()
